# Insert a new "Review Data" slide (Title + Content layout) right before the
# final "PRISMA Flow Diagram" slide, pushing the latter down to become the
# new last slide.

$p = $ppt.ActivePresentation

# The presentation currently has 14 slides; the PRISMA Flow Diagram slide is
# last (index 14). Add the new slide at that position so it lands just
# before it (ppLayoutText = 2, i.e. the "Title and Content" layout already
# used by the PRISMA slide).
$newSlide = $p.Slides.Add($p.Slides.Count, 2)

# Title placeholder.
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Review Data"

# Body / content placeholder - five bullet paragraphs.
$bullets = @(
    "4392 articles are identified.",
    "864 duplicates are detected.",
    "510 are deleted.",
    "354 are resolved.",
    "3882 articles to screen."
)
$newSlide.Shapes.Item(2).TextFrame.TextRange.Text = [string]::Join("`r", $bullets)
